$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COVID Resources-HCP")

# --- Row 100: Self-Compassion Guided Meditations and Exercises (Kristen Neff) ---
$ws.Range("A100").Value = "Health Care Assistants"
$ws.Range("B100").Value = "Healthcare Provider Wellness"
$ws.Range("C100").Value = "United States"
$ws.Range("D100").Value = "Self-Compassion Guided Meditations and Exercises"
$ws.Range("E100").Value = "Kristen Neff"
$ws.Range("F100").Value = "Education"
$ws.Hyperlinks.Add($ws.Range("G100"), "https://self-compassion.org/category/exercises/", "", "", "https://self-compassion.org/category/exercises/")
$ws.Rows.Item(100).RowHeight = 32

# --- Row 101: How to Make Stress Your Friend (Kelly McGonigal) ---
$ws.Range("A101").Value = "Health Care Assistants"
$ws.Range("B101").Value = "Healthcare Provider Wellness"
$ws.Range("C101").Value = "International"
$ws.Range("D101").Value = "How to Make Stress Your Friend"
$ws.Range("E101").Value = "Kelly McGonigal"
$ws.Range("F101").Value = "Education"
$ws.Hyperlinks.Add($ws.Range("G101"), "https://www.ted.com/talks/kelly_mcgonigal_how_to_make_stress_your_friend?language=en", "t-6077", "", "https://www.ted.com/talks/kelly_mcgonigal_how_to_make_stress_your_friend?language=en#t-6077")
$ws.Rows.Item(101).RowHeight = 32

# --- Row 102: Working in Health Care during COVID-19 ---
$ws.Range("A102").Value = "Health Care Assistants"
$ws.Range("B102").Value = "Healthcare Provider Wellness"
$ws.Range("C102").Value = "British Columbia"
$ws.Range("D102").Value = "Working in Health Care during COVID-19"
$ws.Range("E102").Value = "BC Centre for Palliative Care"
$ws.Range("F102").Value = "Education"
$ws.Hyperlinks.Add($ws.Range("G102"), "https://www.youtube.com/watch?v=mAydHEj7JqI&feature=youtu.be", "", "", "https://www.youtube.com/watch?v=mAydHEj7JqI&feature=youtu.be")
$ws.Rows.Item(102).RowHeight = 32

# --- Row 103: HCA COVID Needs Assessment ---
$ws.Range("A103").Value = "Health Care Assistants"
$ws.Range("B103").Value = "Healthcare Provider Wellness"
$ws.Range("C103").Value = "British Columbia"
$ws.Range("D103").Value = "HCA COVID Needs Assessment"
$ws.Range("E103").Value = "BC Centre for Palliative Care"
$ws.Range("F103").Value = "Education"
$hcaUrl = "https://ihsts.sharepoint.com/sites/bc-cpc/Shared%20Documents/Forms/AllItems.aspx?id=%2Fsites%2Fbc%2Dcpc%2FShared%20Documents%2FCommunications%2F2020%20Communications%2FStrategic%20Initiatives%2FHCA%20needs%20assessment%20report%2FHCA%20COVID%20needs%20assessment%2Epdf&parent=%2Fsites%2Fbc%2Dcpc%2FShared%20Documents%2FCommunications%2F2020%20Communications%2FStrategic%20Initiatives%2FHCA%20needs%20assessment%20report&p=true&originalPath=aHR0cHM6Ly9paHN0cy5zaGFyZXBvaW50LmNvbS86Yjovcy9iYy1jcGMvRWRkS0NPZ1hfXzFJanNvV1RHNnl3RDhCZGF6RjZBRlhFdUszUjVhUGpQaUZUUT9ydGltZT1WaHB6TmQ4ejJFZw"
$ws.Hyperlinks.Add($ws.Range("G103"), $hcaUrl, "", "", $hcaUrl)
$ws.Rows.Item(103).RowHeight = 192

# --- Update frozen-pane scroll position and selection to match the new bottom of the table ---
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 103
$ws.Range("C104:G108").Select()
$ws.Range("C104").Activate()
